$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Easy Driving Mode" team member list (K3:K7) - fill in the remaining members
$ws.Range("K5").Value = "Youssef Yasser"
$ws.Range("K6").Value = "Mohamed Aladdin"
$ws.Range("K7").Value = "Youssef Moro"

# "obstacle avoidance" (Ultrasonic sensor) team member list (K11:K15)
$ws.Range("K13").Value = "Youssef Yasser"

# Leave the selection on K10 and scrolled back to the top of the sheet
$ws.Range("K10").Select() | Out-Null
